$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts old D..U (months) right to E..V
# and keeps A (S.No), B (Client Name), C (PO No) in place.
$ws.Columns("D").Insert()

# New "Project Owner" header + filter values (pivot-table style "filter" column)
$ws.Range("D1").Value = "Project Owner"
$ws.Range("D2").Value = "ProjectOwner2"
$ws.Range("D3").Value = "ProjectOwner2"
$ws.Range("D4").Value = "ProjectOwner1"

# A handful of monthly forecast cells were also re-bucketed as part of this
# change (not just shifted one column over) - apply those explicit overrides
# on top of the shifted grid.
$ws.Range("G2").Value = 0

$ws.Range("F3").Value = 0
$ws.Range("I3").Value = 10000
$ws.Range("K3").Value = 15000
$ws.Range("M3").Value = 10000
